# "#2 ordering rows in master file by county and station"
#
# The underlying data change is: the "maxt"/"mint" temperature columns
# (and their legends) are renamed to "maxtp"/"mintp" (one character
# longer each). Everything else about the row/column layout stays the
# same; only the header text, the legend sheet text, the best-fit
# column widths for those two columns, and the last-used-cell
# selections change.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("dly7404")   # data sheet: county/station/date/...
$ws2 = $wb.Worksheets.Item("Sheet1")    # legend sheet describing each column

# --- Data sheet header row: rename "maxt" -> "maxtp" and "mint" -> "mintp" ---
# (mintp is written first, then maxtp, so new shared-string entries are
#  appended in the same order Excel produced them in the target file)
$ws1.Cells.Item(1, 9).Value = "mintp"
$ws1.Cells.Item(1, 7).Value = "maxtp"

# --- Legend sheet: update the matching descriptive rows ---
$ws2.Cells.Item(9, 1).Value = "mintp: -   Minimum Temperature (C)"
$ws2.Cells.Item(8, 1).Value = "maxtp: -   Maximum Temperature (C)"

# --- Column widths grow by one character now that the headers are longer ---
$ws1.Columns.Item(7).ColumnWidth = 5.25
$ws1.Columns.Item(9).ColumnWidth = 4.8

# --- Restore the selected/active cell on each sheet ---
$ws2.Activate() | Out-Null
$ws2.Range("A8").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("C8").Select() | Out-Null
